$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Page 1")

# Update exam grades (and related columns) per row
$ws.Range("J6").Value = 12
$ws.Range("J9").Value = 12
$ws.Range("J11").Value = 8
$ws.Range("H14").Value = 15
$ws.Range("J14").Value = 8
$ws.Range("J17").Value = 12
$ws.Range("H27").Value = 15
$ws.Range("J27").Value = 14
$ws.Range("H28").Value = 15
$ws.Range("J29").Value = 8
$ws.Range("J30").Value = 12
$ws.Range("J31").Value = 12
$ws.Range("H32").Value = 15
$ws.Range("J32").Value = 14
$ws.Range("J33").Value = 12
$ws.Range("J39").Value = 10
$ws.Range("J41").Value = 16
$ws.Range("J45").Value = 8
$ws.Range("J46").Value = 10
$ws.Range("J48").Value = 14
$ws.Range("J50").Value = 14
$ws.Range("J51").Value = 14
$ws.Range("J52").Value = 6
$ws.Range("J54").Value = 14
$ws.Range("H55").Value = 15
$ws.Range("J55").Value = 12
$ws.Range("J56").Value = 10
$ws.Range("H58").Value = 15
$ws.Range("J58").Value = 12
$ws.Range("I61").Value = 16
$ws.Range("J62").Value = 12
$ws.Range("J63").Value = 12
$ws.Range("J64").Value = 10
$ws.Range("J65").Value = 12
$ws.Range("J67").Value = 12
$ws.Range("J69").Value = 8
$ws.Range("J70").Value = 12
$ws.Range("J73").Value = 16
$ws.Range("J75").Value = 14
$ws.Range("J76").Value = 10
$ws.Range("J77").Value = 12
$ws.Range("J78").Value = 12
$ws.Range("J80").Value = 12
$ws.Range("J82").Value = 10
$ws.Range("J85").Value = 16
$ws.Range("J86").Value = 16
$ws.Range("J88").Value = 14
$ws.Range("J89").Value = 12
$ws.Range("H91").Value = 15
$ws.Range("H92").Value = 15
$ws.Range("J92").Value = 14
$ws.Range("J93").Value = 10
$ws.Range("J94").Value = 12
$ws.Range("J96").Value = 6
$ws.Range("H97").Value = 15
$ws.Range("J97").Value = 10
$ws.Range("J98").Value = 10

# Update the active sheet view (frozen pane top-left cell and selection)
$ws.Activate()
$aw = $excel.ActiveWindow
$aw.ScrollRow = 82
$aw.ScrollColumn = 1
$ws.Range("H91").Select()

Write-Host "Applied database exam grade updates"